# Update "想去人数" (F column) values as published on gh-pages output
# generated at 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 20346
$ws1.Range("F4").Value = 805
$ws1.Range("F7").Value = 18
$ws1.Range("F8").Value = 7654
$ws1.Range("F10").Value = 740
$ws1.Range("F11").Value = 281
$ws1.Range("F14").Value = 132
$ws1.Range("F15").Value = 17
$ws1.Range("F17").Value = 200
$ws1.Range("F18").Value = 1344
$ws1.Range("F26").Value = 1125
$ws1.Range("F31").Value = 571
$ws1.Range("F33").Value = 4883
$ws1.Range("F37").Value = 12742
$ws1.Range("F42").Value = 277
$ws1.Range("F44").Value = 4013

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 207

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 20346
$ws4.Range("F4").Value = 805
$ws4.Range("F7").Value = 18
$ws4.Range("F8").Value = 7654
$ws4.Range("F10").Value = 740
$ws4.Range("F11").Value = 281
$ws4.Range("F14").Value = 132
$ws4.Range("F15").Value = 17
$ws4.Range("F17").Value = 200
$ws4.Range("F18").Value = 1344
$ws4.Range("F26").Value = 1125
$ws4.Range("F30").Value = 207
$ws4.Range("F32").Value = 571
$ws4.Range("F36").Value = 4883
$ws4.Range("F40").Value = 12742
$ws4.Range("F45").Value = 277
$ws4.Range("F47").Value = 4013
